$d = $word.ActiveDocument

# The page used to end with a blank spacer paragraph, a "Ver no Jupiter ..."
# line and a "© 2020 ..." footer line, right after the
# "LOB1012: Estatística (Requisito fraco)" requirement paragraph. Those three
# paragraphs (the blank spacer + the two text paragraphs) are being removed
# entirely, leaving the requirement paragraph followed directly by the
# pre-existing blank paragraph / page-break paragraph that close the
# document.

$marker = "Ver no Jupiter Salvar em pdf Salvar em docx"
$footer = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$r = $d.Content
[void]$r.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($r.Find.Found) {
    # Pull the start back one character so the preceding blank paragraph's
    # mark is swallowed too.
    [void]$r.MoveStart(1, -1)
    # Extend the end past this paragraph's own mark, the whole
    # copyright/footer paragraph text, and that paragraph's trailing mark.
    [void]$r.MoveEnd(1, 1 + $footer.Length + 1)

    $victim = $d.Range($r.Start, $r.End)
    $victim.Delete()
}
